# ============================================================================
# feat: add 2022-Q1 data
#
# The previous "总计" (Total) sheet (sheetId=5) is renamed to "2022-Q1" and
# repopulated with that quarter's per-fund holding breakdown. A brand new
# "总计" sheet (sheetId=6) is appended right after it, carrying forward the
# old aggregate rows plus a new leading row summarizing 2022-Q1.
# ============================================================================

$wb = $excel.ActiveWorkbook

# Donor cells that already carry the "bold header / bordered index column"
# style used throughout this workbook (s="2" in the OOXML) - used purely as
# formatting sources via Copy + PasteSpecial(xlPasteFormats = -4122).
$styleSource = $wb.Worksheets.Item("2021-Q4")
$styleHeaderCell = $styleSource.Range("B1")
$styleIndexCell  = $styleSource.Range("A2")
$blankCell       = $styleSource.Range("Z100")

# Force a cell to hold a *text* value even when the string looks numeric
# (e.g. "51.61", "001071") - Excel would otherwise silently coerce it to a
# number (and drop leading zeros). Flip the cell to text format, assign the
# value, then strip the leftover custom number format by pasting the format
# of a pristine, never-touched cell back over it, leaving the cell with no
# style at all - matching the source data (bare <c> with no s="...").
function Set-TextCell($cell, [string]$val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $blankCell.Copy()
    $cell.PasteSpecial(-4122)
}

function Set-HeaderStyle($cell) {
    $styleHeaderCell.Copy()
    $cell.PasteSpecial(-4122)
}

function Set-IndexStyle($cell) {
    $styleIndexCell.Copy()
    $cell.PasteSpecial(-4122)
}

$fundData = @(
    @('001071','华安媒体互联网混合','51.61','92.88','3.69','1.9044',10),
    @('590001','中邮核心优选混合','16.24','65.06','4.88','0.7925',4),
    @('159766','富国中证旅游主题交易型开放式指数证券投资基金','17.47','99.35','3.65','0.6377',10),
    @('590005','中邮核心主题混合','6.86','65.06','4.98','0.3416',3),
    @('161609','融通动力先锋混合','7.12','80.93','4.00','0.2848',1),
    @('210003','金鹰行业优势混合','6.62','88.02','4.30','0.2847',6),
    @('000968','广发中证养老产业指数A','10.39','94.08','1.71','0.1777',3),
    @('001152','融通新区域新经济灵活配置混合','2.81','80.98','4.01','0.1127',1),
    @('004206','华商元亨灵活配置混合','5.64','29.94','1.89','0.1066',2),
    @('770001','德邦优化灵活配置混合','2.49','86.80','3.00','0.0747',7),
    @('562510','华夏中证旅游主题ETF','1.71','98.99','3.63','0.0621',10),
    @('011155','金鹰责任投资混合A','1.23','82.39','4.27','0.0525',8),
    @('003598','华商润丰灵活配置混合A','3.17','37.86','1.63','0.0517',6),
    @('007509','华商润丰灵活配置混合C','3.09','37.86','1.63','0.0504',6),
    @('002681','金鹰元和灵活配置混合A','0.56','81.63','4.70','0.0263',6),
    @('003513','中邮消费升级灵活配置混合','0.41','54.24','4.81','0.0197',5),
    @('011156','金鹰责任投资混合C','0.39','82.39','4.27','0.0167',8),
    @('002982','广发中证养老产业指数C','0.88','94.08','1.71','0.0150',3),
    @('516560','华宝养老ETF','0.75','97.92','1.78','0.0134',3),
    @('002682','金鹰元和灵活配置混合C','0.25','81.63','4.70','0.0118',6),
    @('001448','华商双翼平衡混合','0.38','39.74','3.04','0.0116',1),
    @('001664','平安鑫安混合A','0.86','29.46','0.92','0.0079',9),
    @('510190','华安上证龙头ETF','0.53','97.53','1.28','0.0068',5),
    @('007049','平安鑫安混合E','0.50','29.46','0.92','0.0046',9),
    @('673120','西部利得新富灵活配置混合','0.02','80.77','3.44','0.0007',10),
    @('001665','平安鑫安混合C','0.02','29.46','0.92','0.0002',9)
)

$totalData = @(
    @('2022-Q1', 26, 5.07),
    @('2021-Q4', 5, 2.64),
    @('2021-Q3', 3, 2.59),
    @('2021-Q2', 2, 2.68),
    @('2021-Q1', 11, 4.1)
)

# ----------------------------------------------------------------------
# Step 1: the existing "总计" sheet becomes "2022-Q1" (same sheetId=5,
# same position), repopulated with the per-fund breakdown for that quarter.
# ----------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

$headers1 = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers1.Count; $c++) {
    $cell = $q1.Cells.Item(1, 2 + $c)
    $cell.Value = $headers1[$c]
    Set-HeaderStyle $cell
}

for ($i = 0; $i -lt $fundData.Count; $i++) {
    $r = $i + 2
    $row = $fundData[$i]

    $idxCell = $q1.Cells.Item($r, 1)
    $idxCell.Value = $i
    Set-IndexStyle $idxCell

    Set-TextCell $q1.Cells.Item($r, 2) $row[0]
    Set-TextCell $q1.Cells.Item($r, 3) $row[1]
    Set-TextCell $q1.Cells.Item($r, 4) $row[2]
    Set-TextCell $q1.Cells.Item($r, 5) $row[3]
    Set-TextCell $q1.Cells.Item($r, 6) $row[4]
    Set-TextCell $q1.Cells.Item($r, 7) $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
}

# ----------------------------------------------------------------------
# Step 2: append a brand new "总计" sheet right after "2022-Q1" holding the
# aggregate rows (old ones shifted down, new 2022-Q1 row on top). Duplicate
# "2022-Q1" itself (via Copy) so sheet-level properties - sheetPr,
# sheetFormatPr, pageMargins - come along for free, then wipe + refill it.
# ----------------------------------------------------------------------
$q1.Copy($null, $q1)
$total = $wb.Worksheets.Item($q1.Index + 1)
$total.Name = "总计"
$total.Cells.Clear()

$headers2 = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($c = 0; $c -lt $headers2.Count; $c++) {
    $cell = $total.Cells.Item(1, 2 + $c)
    $cell.Value = $headers2[$c]
    Set-HeaderStyle $cell
}

for ($i = 0; $i -lt $totalData.Count; $i++) {
    $r = $i + 2
    $row = $totalData[$i]

    $idxCell = $total.Cells.Item($r, 1)
    $idxCell.Value = $i
    Set-IndexStyle $idxCell

    Set-TextCell $total.Cells.Item($r, 2) $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
}

# Restore the originally active sheet (index 1) so the new/copied sheets
# don't linger as the tab-selected / active one.
$wb.Worksheets.Item(1).Activate()

Write-Output "completed"
